$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")  # "Artfynd" is the (only) sheet in this workbook

# Row 12
$ws.Range("A12").Value = 112128498
$ws.Range("B12").Value = 90018
$ws.Range("C12").Value = "Ovaliderad"
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 1339
$ws.Range("F12").Value = "Brandticka"
$ws.Range("G12").Value = "Pycnoporellus fulgens"
$ws.Range("H12").Value = "(Fr.) Donk"
$ws.Range("P12").Value = "Skansberg, Ö om, Srm"
$ws.Range("Q12").Value = 657133.965947984
$ws.Range("R12").Value = 6571270.665650261
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = "Stockholm"
$ws.Range("U12").Value = "Botkyrka"
$ws.Range("V12").Value = "Södermanland"
$ws.Range("W12").Value = "Botkyrka"
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = "2023-08-28"
$ws.Range("Z12").NumberFormat = "@"
$ws.Range("Z12").Value = "00:00"
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = "2023-08-28"
$ws.Range("AB12").NumberFormat = "@"
$ws.Range("AB12").Value = "00:00"
$ws.Range("AC12").Value = "På granlåga. En del årsfärska dödade granar av granbarkborre. Gott om död ved i form av torrträd och lågor av gran."
$ws.Range("AD12").Value = $False
$ws.Range("AE12").Value = $False
$ws.Range("AG12").Value = $False
$ws.Range("AW12").Value = "Bo Törnquist"
$ws.Range("AX12").Value = "Bo Törnquist"

# Row 13
$ws.Range("A13").Value = 112128627
$ws.Range("B13").Value = 90087
$ws.Range("C13").Value = "Ovaliderad"
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 3298
$ws.Range("F13").Value = "Trådticka"
$ws.Range("G13").Value = "Climacocystis borealis"
$ws.Range("H13").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("P13").Value = "Skansberg, Ö om, Srm"
$ws.Range("Q13").Value = 657181.8354338486
$ws.Range("R13").Value = 6571192.256437058
$ws.Range("S13").Value = 5
$ws.Range("T13").Value = "Stockholm"
$ws.Range("U13").Value = "Botkyrka"
$ws.Range("V13").Value = "Södermanland"
$ws.Range("W13").Value = "Botkyrka"
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2023-08-28"
$ws.Range("Z13").NumberFormat = "@"
$ws.Range("Z13").Value = "00:00"
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = "2023-08-28"
$ws.Range("AB13").NumberFormat = "@"
$ws.Range("AB13").Value = "00:00"
$ws.Range("AC13").Value = "På nedre delen av torrgran."
$ws.Range("AD13").Value = $False
$ws.Range("AE13").Value = $False
$ws.Range("AG13").Value = $False
$ws.Range("AW13").Value = "Bo Törnquist"
$ws.Range("AX13").Value = "Bo Törnquist"

# Row 14
$ws.Range("A14").Value = 112128551
$ws.Range("B14").Value = 90666
$ws.Range("C14").Value = "Ovaliderad"
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 4364
$ws.Range("F14").Value = "Dropptaggsvamp"
$ws.Range("G14").Value = "Hydnellum ferrugineum"
$ws.Range("H14").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P14").Value = "Skansberg, Ö om, Srm"
$ws.Range("Q14").Value = 657162.1501835568
$ws.Range("R14").Value = 6571271.319579108
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = "Stockholm"
$ws.Range("U14").Value = "Botkyrka"
$ws.Range("V14").Value = "Södermanland"
$ws.Range("W14").Value = "Botkyrka"
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = "2023-08-28"
$ws.Range("Z14").NumberFormat = "@"
$ws.Range("Z14").Value = "00:00"
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = "2023-08-28"
$ws.Range("AB14").NumberFormat = "@"
$ws.Range("AB14").Value = "00:00"
$ws.Range("AD14").Value = $False
$ws.Range("AE14").Value = $False
$ws.Range("AG14").Value = $False
$ws.Range("AW14").Value = "Bo Törnquist"
$ws.Range("AX14").Value = "Bo Törnquist"

# Row 15
$ws.Range("A15").Value = 112128530
$ws.Range("B15").Value = 89802
$ws.Range("C15").Value = "Ovaliderad"
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 5420
$ws.Range("F15").Value = "Grovticka"
$ws.Range("G15").Value = "Phaeolus schweinitzii"
$ws.Range("H15").Value = "(Fr.) Pat."
$ws.Range("P15").Value = "Skansberg, Ö om, Srm"
$ws.Range("Q15").Value = 657143.932178961
$ws.Range("R15").Value = 6571277.734310649
$ws.Range("S15").Value = 5
$ws.Range("T15").Value = "Stockholm"
$ws.Range("U15").Value = "Botkyrka"
$ws.Range("V15").Value = "Södermanland"
$ws.Range("W15").Value = "Botkyrka"
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = "2023-08-28"
$ws.Range("Z15").NumberFormat = "@"
$ws.Range("Z15").Value = "00:00"
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = "2023-08-28"
$ws.Range("AB15").NumberFormat = "@"
$ws.Range("AB15").Value = "00:00"
$ws.Range("AC15").Value = "På högstubbe av tall."
$ws.Range("AD15").Value = $False
$ws.Range("AE15").Value = $False
$ws.Range("AG15").Value = $False
$ws.Range("AW15").Value = "Bo Törnquist"
$ws.Range("AX15").Value = "Bo Törnquist"

# Row 16
$ws.Range("A16").Value = 112128573
$ws.Range("B16").Value = 90668
$ws.Range("C16").Value = "Ovaliderad"
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 788
$ws.Range("F16").Value = "Gul taggsvamp"
$ws.Range("G16").Value = "Hydnellum geogenium"
$ws.Range("H16").Value = "(Fr.) Banker"
$ws.Range("P16").Value = "Skansberg, Ö om, Srm"
$ws.Range("Q16").Value = 657133.5412061054
$ws.Range("R16").Value = 6571218.932669931
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = "Stockholm"
$ws.Range("U16").Value = "Botkyrka"
$ws.Range("V16").Value = "Södermanland"
$ws.Range("W16").Value = "Botkyrka"
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2023-08-28"
$ws.Range("Z16").NumberFormat = "@"
$ws.Range("Z16").Value = "00:00"
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2023-08-28"
$ws.Range("AB16").NumberFormat = "@"
$ws.Range("AB16").Value = "00:00"
$ws.Range("AD16").Value = $False
$ws.Range("AE16").Value = $False
$ws.Range("AG16").Value = $False
$ws.Range("AW16").Value = "Bo Törnquist"
$ws.Range("AX16").Value = "Bo Törnquist"

# Row 17
$ws.Range("A17").Value = 112128708
$ws.Range("B17").Value = 90687
$ws.Range("C17").Value = "Ovaliderad"
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 5964
$ws.Range("F17").Value = "Fjällig taggsvamp s.str."
$ws.Range("G17").Value = "Sarcodon imbricatus s.str."
$ws.Range("H17").Value = "(L.:Fr.) P.Karst."
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "1"
$ws.Range("J17").Value = "fruktkroppar"
$ws.Range("P17").Value = "Skansberg, Ö om, Srm"
$ws.Range("Q17").Value = 657216.3353243669
$ws.Range("R17").Value = 6571312.987947896
$ws.Range("S17").Value = 5
$ws.Range("T17").Value = "Stockholm"
$ws.Range("U17").Value = "Botkyrka"
$ws.Range("V17").Value = "Södermanland"
$ws.Range("W17").Value = "Botkyrka"
$ws.Range("Y17").NumberFormat = "@"
$ws.Range("Y17").Value = "2023-08-28"
$ws.Range("Z17").NumberFormat = "@"
$ws.Range("Z17").Value = "00:00"
$ws.Range("AA17").NumberFormat = "@"
$ws.Range("AA17").Value = "2023-08-28"
$ws.Range("AB17").NumberFormat = "@"
$ws.Range("AB17").Value = "00:00"
$ws.Range("AD17").Value = $False
$ws.Range("AE17").Value = $False
$ws.Range("AG17").Value = $False
$ws.Range("AW17").Value = "Bo Törnquist"
$ws.Range("AX17").Value = "Bo Törnquist"
